$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("A2").Value = "MARQUEZ GONZALEZ CARLOS DANIEL"
$ws.Range("B2").Value = "BERA OUTLET"
$ws.Range("D2").Value = "Licuadora"

# Remove row 3 entirely (shift cells up)
$ws.Rows("3:3").Delete()
